# Commit: "setting up table for collecting Q factor data"
# Adds explanatory text rows and a new results table (for Q factor data)
# below the existing FWHM table on Sheet1.
#
# NOTE: cell writes below are deliberately ordered to reproduce the exact
# shared-string allocation order of the original edit (new strings are
# appended to sharedStrings.xml in the order the cells are first written).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Explanatory paragraphs (rows 64-69), all in column A ---
$ws.Range("A64").Value = "The analysis results below will calculate fsr and fwhm data from the saved wavelength sweep data, in a similar way to above. "
$ws.Range("A65").Value = "The difference however, is that to find the peaks and their FWHM, I have converted powers first into mW from dBm."
$ws.Range("A66").Value = "The rest is as above, and can refer to the notebook saved in the same directory as the one used above, but is this time called: data_analysis_notebook_Q_factor."
$ws.Range("A67").Value = "In addition to calculating the above. The notebook finds a list of Q factors by taking each peak_wavelengths/peak FWHM. And finds the mean Q factor and standard error in the Q factor in the same way the mean and standard error for the fsr are found from the list of fsrs."
$ws.Range("A68").Value = "Also, when I input the ring radius, and its error. The notebook calculates Q^3/R^2 and its associated uncertainty as well."
$ws.Range("A69").Value = "Also, note this time, I'll choose the prominence and it will be in mW."

# --- New table header row (row 72), first part (A72:G72) ---
$ws.Range("A72").Value = "Data CSV Filename"
$ws.Range("B72").Value = "Wavelength step size/nm"
$ws.Range("C72").Value = "Start array index"
$ws.Range("D72").Value = "End array index"
$ws.Range("E72").Value = "Start wavelength/nm"
$ws.Range("F72").Value = "End wavelength/nm"
$ws.Range("G72").Value = "prominence/mW"

# --- Row 70: explanatory note about choosing the prominence ---
$ws.Range("A70").Value = "I will choose the prominence roughly by looking at height span of roughly biggest height span noise. But it is very rough, and main check is whether code visually appears to find peaks."

# Row 71 is intentionally left blank (matches the gap seen in the diff).

# --- Rest of the new table header row (row 72) ---
$ws.Range("H72").Value = "distance"
$ws.Range("I72").Value = "approx_fsr/nm"
$ws.Range("J72").Value = "fsr_mean/nm"
$ws.Range("K72").Value = "fsr_std error/nm"
$ws.Range("L72").Value = "double count check passed?"
$ws.Range("M72").Value = "mean FWHM/nm"
$ws.Range("N72").Value = "FWHM error/nm"
$ws.Range("O72").Value = "Q"
$ws.Range("P72").Value = "Q error"
$ws.Range("Q72").Value = "Q^3/R^2 (micrometres^-2)"
$ws.Range("R72").Value = "Q^3/R^2 error (micrometres^-2)"

# --- Row 73: reuses the existing "(approx_fsr/2)/wavelength step size" string in H ---
$ws.Range("H73").Value = "(approx_fsr/2)/wavelength step size"

# Leave the final selection on R72, matching the author's saved cursor position.
$ws.Range("R72").Select()
